# "neue funktion um einzelnen text anzuzeigen"
# Replace the second Q&A card (image-based getter/setter card) and the
# following two placeholder cards (Frage2/Antwort2, Frage3/Antwort3) with
# three new simple text cards (Antwort0/1/2) plus a fresh free-text card.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$question = "Create getter and setter methods for the class."
$bigQuestion = "Wie kann man den Ding am besten machen ?  Wie geht der blabla einfacher?  Was ist die Antwort auf alles?  "

# Row 2: keep the question text in A2, drop the image paths, add Antwort0.
$ws.Range("A2").Value = $question
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = "Antwort0"
$ws.Range("D2").ClearContents()

# Row 3: repeat the same question with Antwort1 (was the image/placeholder card).
$ws.Range("A3").Value = $question
$ws.Range("B3").ClearContents()
$ws.Range("C3").Value = "Antwort1"
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()

# Row 4: repeat the same question again with Antwort2 (was Frage2/Antwort2).
$ws.Range("A4").Value = $question
$ws.Range("C4").Value = "Antwort2"

# Row 5: new free-text question/answer pair (was Frage3/Antwort3).
$ws.Range("A5").Value = $bigQuestion
$ws.Range("C5").Value = "s"

# The question column in rows 3-5 now carries the same left/top wrapped
# formatting as the question cell above it (row 2), matching row 2's style.
$ws.Range("A3").HorizontalAlignment = -4131
$ws.Range("A3").VerticalAlignment = -4160
$ws.Range("A4").HorizontalAlignment = -4131
$ws.Range("A4").VerticalAlignment = -4160
$ws.Range("A5").HorizontalAlignment = -4131
$ws.Range("A5").VerticalAlignment = -4160

# Row heights follow the card's question length, like the other cards
# already on the sheet (row 2 and the newly duplicated rows 3/4 use the
# short-question height, row 5 takes the tall free-text height).
$ws.Rows.Item(2).RowHeight = 28.35
$ws.Rows.Item(3).RowHeight = 28.35
$ws.Rows.Item(4).RowHeight = 28.35
$ws.Rows.Item(5).RowHeight = 55.2

# Move the selection to A5, matching the new active cell in the sheet view.
$ws.Range("A5").Select()
